$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row 4 (old row 4 -> row 5, old row 5 -> row 6)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Row 1 header: new text, merged A1:I1, new style
# ---------------------------------------------------------------------------
$r1 = $ws.Range("A1:I1")
$r1.UnMerge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in C. Tbilisi Municipality"
$r1.Merge()
$ws.Rows.Item(1).RowHeight = 51
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.Interior.Pattern = -4142
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4108
$r1.WrapText = $true
$r1.Borders.Item(7).LineStyle = -4142
$r1.Borders.Item(8).LineStyle = -4142
$r1.Borders.Item(9).LineStyle = -4142
$r1.Borders.Item(10).LineStyle = -4142

# ---------------------------------------------------------------------------
# 3. Row 2 "(End of year, persons)" keeps its formatting already - no change
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4. A3 (blank header cell above years): font becomes Sylfaen 11
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. New row 4: "family with disabilities Persons " + values
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$a4 = $ws.Range("A4")
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.Color = 0
$a4.Interior.Pattern = 1
$a4.Interior.ThemeColor = 1
$a4.Interior.TintAndShade = 0
$a4.Borders.Item(7).LineStyle = -4142
$a4.Borders.Item(8).LineStyle = 1
$a4.Borders.Item(8).Weight = 2
$a4.Borders.Item(9).LineStyle = -4142
$a4.Borders.Item(10).LineStyle = -4142
$a4.HorizontalAlignment = -4131
$a4.VerticalAlignment = -4108
$a4.WrapText = $true
$ws.Rows.Item(4).RowHeight = 24.75

$row4vals = @(11238,10920,10728,11337,11430,11595,11576,11990)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $row4vals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Color = 0
    $cell.Interior.Pattern = 1
    $cell.Interior.ThemeColor = 1
    $cell.Interior.TintAndShade = 0
    $cell.NumberFormat = "#\ ##0"
    $cell.Borders.Item(7).LineStyle = -4142
    $cell.Borders.Item(8).LineStyle = -4142
    $cell.Borders.Item(9).LineStyle = -4142
    $cell.Borders.Item(10).LineStyle = -4142
    $cell.HorizontalAlignment = 1
    $cell.VerticalAlignment = -4107
    $cell.WrapText = $false
}

# ---------------------------------------------------------------------------
# 6. Row 5 (old row 4): "disabilities Persons " + new values
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$a5 = $ws.Range("A5")
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.Color = 0
$a5.Borders.Item(7).LineStyle = -4142
$a5.Borders.Item(8).LineStyle = -4142
$a5.Borders.Item(9).LineStyle = 1
$a5.Borders.Item(9).Weight = 2
$a5.Borders.Item(10).LineStyle = -4142
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$ws.Rows.Item(5).RowHeight = 21

$row5vals = @(12448,12095,11864,12497,12565,12726,12716)
$cols2 = @("B","C","D","E","F","G","H")
for ($i = 0; $i -lt 7; $i++) {
    $cell = $ws.Range($cols2[$i] + "5")
    $cell.Value = $row5vals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Color = 0
    $cell.Borders.Item(7).LineStyle = -4142
    $cell.Borders.Item(8).LineStyle = -4142
    $cell.Borders.Item(9).LineStyle = -4142
    $cell.Borders.Item(10).LineStyle = -4142
    $cell.NumberFormat = "#\ ##0"
    $cell.HorizontalAlignment = 1
    $cell.VerticalAlignment = -4107
    $cell.WrapText = $false
}
$i5 = $ws.Range("I5")
$i5.Value = 13143
$i5.Font.Name = "Arial"
$i5.Font.Size = 10
$i5.Font.Color = 0
$i5.Borders.Item(7).LineStyle = -4142
$i5.Borders.Item(8).LineStyle = -4142
$i5.Borders.Item(9).LineStyle = 1
$i5.Borders.Item(9).Weight = 2
$i5.Borders.Item(10).LineStyle = -4142
$i5.NumberFormat = "#\ ##0"
$i5.HorizontalAlignment = 1
$i5.VerticalAlignment = -4107
$i5.WrapText = $false

# ---------------------------------------------------------------------------
# 7. Row 6 (old row 5): Source row - text unchanged, formatting tweaks
# ---------------------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Borders.Item(8).LineStyle = -4142
$ws.Rows.Item(6).RowHeight = 27.75

Write-Host "done"
